$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate "vamc-upgrade Sprint 8" three times to create Sprint 9, 10, 11
#    Each copy is inserted immediately before "DevTeam" (re-fetched fresh
#    each time since the worksheet collection shifts after every Copy).
# ---------------------------------------------------------------------------
$sprint8 = $wb.Worksheets.Item("vamc-upgrade Sprint 8")

$devteam = $wb.Worksheets.Item("DevTeam")
$sprint8.Copy($devteam, $null)
$sprint9 = $wb.Worksheets.Item(2)
$sprint9.Name = "vamc-upgrade Sprint 9"

$devteam = $wb.Worksheets.Item("DevTeam")
$sprint8.Copy($devteam, $null)
$sprint10 = $wb.Worksheets.Item(3)
$sprint10.Name = "vamc-upgrade Sprint 10"

$devteam = $wb.Worksheets.Item("DevTeam")
$sprint8.Copy($devteam, $null)
$sprint11 = $wb.Worksheets.Item(4)
$sprint11.Name = "vamc-upgrade Sprint 11"

# ---------------------------------------------------------------------------
# 2. Update the capacity numbers on the new sheets. Columns D (Days
#    Available) and F (Available Hours) are formulas and recompute on
#    their own -- only the raw inputs in B (Number of Days in Sprint) and
#    C (Days Off) need to be written.
# ---------------------------------------------------------------------------

# vamc-upgrade Sprint 9 (same sprint length as Sprint 8 -- 9 days)
$sprint9.Range("C2").Value = 0
$sprint9.Range("C3").Value = 1
$sprint9.Range("C4").Value = 9
$sprint9.Range("C5").Value = 1
$sprint9.Range("C6").Value = 0
$sprint9.Range("C7").Value = 0

# vamc-upgrade Sprint 10 (10-day sprint)
$sprint10.Range("B2:B7").Value = 10
$sprint10.Range("C2").Value = 1
$sprint10.Range("C3").Value = 0
$sprint10.Range("C4").Value = 9
$sprint10.Range("C5").Value = 0
$sprint10.Range("C6").Value = 0
$sprint10.Range("C7").Value = 4

# vamc-upgrade Sprint 11 (8-day sprint)
$sprint11.Range("B2:B7").Value = 8
$sprint11.Range("C2").Value = 0.5
$sprint11.Range("C3").Value = 0
$sprint11.Range("C4").Value = 0
$sprint11.Range("C5").Value = 2
$sprint11.Range("C6").Value = 0
$sprint11.Range("C7").Value = 3

# ---------------------------------------------------------------------------
# 3. Add the "sprint date range" note in row 13 of every sprint sheet.
#    Shared strings are created in this exact order so they land at the
#    same indices as the target workbook (34..37).
# ---------------------------------------------------------------------------
$sprint9.Range("B13").Value = "November 20th - Dec 3"
$sprint8.Range("B13").Value = "November 6th - 19th"
$sprint10.Range("B13").Value = "December 4th-17th"
$sprint11.Range("B13").Value = "December 18th-31st"

# ---------------------------------------------------------------------------
# 4. Re-create the hidden _FilterDatabase defined names for the new sheets
#    (AutoFilter is copied onto each new sheet already, but the defined
#    name bookkeeping is not, so it is added back explicitly).
# ---------------------------------------------------------------------------
$sprint9.Names.Add("_xlnm._FilterDatabase", "='vamc-upgrade Sprint 9'!`$A`$1:`$F`$1")
$sprint10.Names.Add("_xlnm._FilterDatabase", "='vamc-upgrade Sprint 10'!`$A`$1:`$F`$1")
$sprint11.Names.Add("_xlnm._FilterDatabase", "='vamc-upgrade Sprint 11'!`$A`$1:`$F`$1")

# ---------------------------------------------------------------------------
# 5. Sheet-view bookkeeping: selection per sheet + which tab is active.
#    Activating a sheet also flips its sheetView's tabSelected flag and the
#    workbook's activeTab, so the order here matters -- Sprint 11 ends up
#    being the active tab, matching the target.
# ---------------------------------------------------------------------------
$sprint8.Activate()
$sprint8.Range("B13").Select()

$sprint9.Activate()
$sprint9.Range("C9").Select()

$sprint10.Activate()
$sprint10.Range("C3").Select()

$sprint11.Activate()
$sprint11.Range("C3").Select()
